# Apply updated cryptocurrency market data (prices & 1h volume %) scraped Fri Jun  9 19:55:00 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.438.73"
$ws.Range("E2").Value = "  -0.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.834.55"
$ws.Range("E3").Value = "  -0.64%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "259.99"
$ws.Range("E5").Value = "  -1.58%  "

$ws.Range("E6").Value = "  +0.19%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5325"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3005"
$ws.Range("E8").Value = "  -6.70%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06856"
$ws.Range("E9").Value = "  +0.74%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.48"
$ws.Range("E10").Value = "  -7.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7341"
$ws.Range("E11").Value = "  -5.75%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.839.39"
$ws.Range("E12").Value = "  -0.22%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07265"
$ws.Range("E13").Value = "  -6.46%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.81"
$ws.Range("E14").Value = "  +0.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.951"
$ws.Range("E15").Value = "  -1.44%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  +0.32%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.86"
$ws.Range("E17").Value = "  -0.82%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  +0.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007861"
$ws.Range("E19").Value = "  -1.36%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.467.30"
$ws.Range("E20").Value = "  -0.31%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.079.57"
$ws.Range("E21").Value = "  +0.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.572"
$ws.Range("E22").Value = "  -1.31%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.955"
$ws.Range("E23").Value = "  -0.92%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.221"
$ws.Range("E24").Value = "  -2.66%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.65"
$ws.Range("E25").Value = "  -0.15%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.191"
$ws.Range("E26").Value = "  +0.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.682"
$ws.Range("E27").Value = "  +0.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.90"
$ws.Range("E28").Value = "  -0.75%  "

$ws.Range("E29").Value = "  -1.25%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.212"
$ws.Range("E30").Value = "  +0.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08799"
$ws.Range("E31").Value = "  +0.70%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.995"
$ws.Range("E32").Value = "  -3.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04780"
$ws.Range("E33").Value = "  -1.16%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.931"
$ws.Range("E34").Value = "  +2.18%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7284"
$ws.Range("E35").Value = "  +1.31%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.126"
$ws.Range("E36").Value = "  -0.34%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.087"
$ws.Range("E37").Value = "  -0.78%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.286"
$ws.Range("E38").Value = "  +2.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01700"
$ws.Range("E39").Value = "  -5.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4700"
$ws.Range("E40").Value = "  -3.50%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9029"
$ws.Range("E41").Value = "  +0.55%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "107.27"
$ws.Range("E42").Value = "  -3.06%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.875"
$ws.Range("E43").Value = "  -2.41%  "

$ws.Range("E44").Value = "  +0.17%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.334"
$ws.Range("E45").Value = "  -3.98%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.950"
$ws.Range("E46").Value = "  -0.66%  "

$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1226"
$ws.Range("E47").Value = "  -0.08%  "

$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4047"
$ws.Range("E48").Value = "  -3.81%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05794"
$ws.Range("E49").Value = "  -1.59%  "

$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.70"
$ws.Range("E50").Value = "  -0.91%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8903"
$ws.Range("E51").Value = "  +0.10%  "
